$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new header cell "Save" in H1, copying the header formatting from G1
# (bold font, centered alignment, thin border) via copy/paste-special so the
# cell reuses the existing style definition exactly.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the "Save" column values for rows 2-7
$values = @(1, 0, 0, 0, 1, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
